$d = $word.ActiveDocument

# Locate the M2Doc field that looks like "{ m:     }" (a Word field whose
# instrText spells out the missing-expression token " m:     ").
$target = $null
foreach ($fld in $d.Fields) {
    if ($fld.Code.Text -match "m:") {
        $target = $fld
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the target field containing 'm:'"
}

# The field begin mark sits one position before the instrText (Code) range.
$fieldRange = $target.Code
$start = $fieldRange.Start - 1

# Remove the whole field (begin mark, instrText runs, end mark) in one shot.
$target.Delete()

# Replace it with plain literal text "{m:}" - the rewriter now emits the
# token delimiters and the expression text as ordinary runs instead of a
# Word field.
$ins = $d.Range($start, $start)
$ins.InsertBefore("{m:}")

# Split "{m:}" into three distinct runs -> "{", "m", ":}" (matching the
# token splitting performed by TokenIteratorFieldRewriterSplit) by
# temporarily bookmarking the split points, then removing the bookmarks.
$b1 = $d.Bookmarks.Add("TempTokenSplit1", $d.Range($start + 1, $start + 1))
$b2 = $d.Bookmarks.Add("TempTokenSplit2", $d.Range($start + 2, $start + 2))
$d.Bookmarks("TempTokenSplit1").Delete()
$d.Bookmarks("TempTokenSplit2").Delete()
